# Insert a new weekly price record as row 152 in the "Papa" (potato) price
# sheet. This pushes the existing rows 152-178 down to 153-179 (dimension
# grows from A1:R178 to A1:R179) and populates the newly-opened row 152
# with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 152..178 down to 153..179, opening up a blank row 152.
$ws.Rows(152).Insert()

# Populate the new row 152 with the new record.
$ws.Cells.Item(152, 1).Value  = 1
$ws.Cells.Item(152, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(152, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(152, 4).Value  = 45005
$ws.Cells.Item(152, 5).Value  = 15
$ws.Cells.Item(152, 6).Value  = 100114001
$ws.Cells.Item(152, 7).Value  = "Papa"
$ws.Cells.Item(152, 8).Value  = "Red Lady"
$ws.Cells.Item(152, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(152, 10).Value = 1250
$ws.Cells.Item(152, 11).Value = 13000
$ws.Cells.Item(152, 12).Value = 14000
$ws.Cells.Item(152, 13).Value = 13600
$ws.Cells.Item(152, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(152, 15).Value = "Región del Maule"
$ws.Cells.Item(152, 16).Value = 544
$ws.Cells.Item(152, 17).Value = 25
$ws.Cells.Item(152, 18).Value = "Hortaliza"
